$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Sheet1"

$ws.Range("A1").Value = "Part Number"
$ws.Range("C1").Value = "Manufacturer"
$ws.Range("D1").Value = "Cost"
$ws.Range("A2").Value = "GPA18GT2060-A-H5"
$ws.Range("B1").Value = "Item Name"
$ws.Range("C2").Value = "MISUMI"
$ws.Range("A3").Value = "GPA60GT2060-B-H8"
$ws.Range("B2").Value = "GT2 Timing Belt Pulley 18 Tooth"
$ws.Range("B3").Value = "GT2 Timing Belt Pulley 20 Tooth"

$ws.Range("D2").Value = 18.1
$ws.Range("C3").Value = "MISUMI"
$ws.Range("D3").Value = 32.18

$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("D2:D3").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

$ws.Columns("A:C").AutoFit()

$ws.Range("C9").Select()
